$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# "partb" starts life as an exact copy of "parta" (same column widths,
# styles, and F/D/X/M/W pipeline-diagram layout), placed right after it.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "partb"

# --- adjust cell contents so "partb" reflects the part-b instruction
#     sequence/timing (wider diagram, extra lw instructions) ---
$ws2.Cells.Item(1, 24).Value = "m"
$ws2.Cells.Item(1, 25).Value = "m"
$ws2.Cells.Item(1, 26).Value = "m"
$ws2.Cells.Item(1, 27).Value = "m"
$ws2.Cells.Item(1, 29).Value = "m"
$ws2.Cells.Item(1, 30).Value = "m"
$ws2.Cells.Item(1, 31).Value = "m"
$ws2.Cells.Item(1, 32).Value = "m"
$ws2.Cells.Item(1, 33).Value = "u"
$ws2.Cells.Item(1, 35).Value = "u"
$ws2.Cells.Item(1, 36).Value = "u"
$ws2.Cells.Item(1, 37).Value = "u"
$ws2.Cells.Item(1, 38).Value = "u"
$ws2.Cells.Item(1, 39).Value = "u"
$ws2.Cells.Item(1, 40).Value = "c"
$ws2.Cells.Item(1, 41).Value = "c"
$ws2.Cells.Item(1, 42).Value = "u"
$ws2.Cells.Item(2, 35).Value = 34
$ws2.Cells.Item(2, 36).Value = 35
$ws2.Cells.Item(2, 37).Value = 36
$ws2.Cells.Item(2, 38).Value = 37
$ws2.Cells.Item(2, 39).Value = 38
$ws2.Cells.Item(2, 40).Value = 39
$ws2.Cells.Item(2, 41).Value = 40
$ws2.Cells.Item(2, 42).Value = 41
$ws2.Cells.Item(9, 1).Value = "lw r5, 8(r5)"
$ws2.Cells.Item(8, 1).Value = "lw r4, 4(r4)"
$ws2.Cells.Item(13, 24).Value = "M"
$ws2.Cells.Item(13, 25).Value = "M"
$ws2.Cells.Item(13, 26).Value = "M"
$ws2.Cells.Item(13, 27).Value = "M"
$ws2.Cells.Item(13, 28).Value = "W"
$ws2.Cells.Item(14, 24).Value = "X"
$ws2.Cells.Item(14, 25).Value = "X"
$ws2.Cells.Item(14, 26).Value = "X"
$ws2.Cells.Item(14, 27).Value = "X"
$ws2.Cells.Item(14, 28).Value = "M"
$ws2.Cells.Item(14, 29).Value = "M"
$ws2.Cells.Item(14, 30).Value = "M"
$ws2.Cells.Item(14, 31).Value = "M"
$ws2.Cells.Item(14, 32).Value = "M"
$ws2.Cells.Item(14, 33).Value = "W"
$ws2.Cells.Item(15, 24).Value = "D"
$ws2.Cells.Item(15, 25).Value = "D"
$ws2.Cells.Item(15, 26).Value = "D"
$ws2.Cells.Item(15, 27).Value = "D"
$ws2.Cells.Item(15, 28).Value = "X"
$ws2.Cells.Item(15, 29).Value = "X"
$ws2.Cells.Item(15, 30).Value = "X"
$ws2.Cells.Item(15, 31).Value = "X"
$ws2.Cells.Item(15, 32).Value = "X"
$ws2.Cells.Item(15, 33).Value = "M"
$ws2.Cells.Item(15, 34).Value = "W"
$ws2.Cells.Item(16, 24).Value = "F"
$ws2.Cells.Item(16, 25).Value = "F"
$ws2.Cells.Item(16, 26).Value = "F"
$ws2.Cells.Item(16, 27).Value = "F"
$ws2.Cells.Item(16, 28).Value = "D"
$ws2.Cells.Item(16, 29).Value = "D"
$ws2.Cells.Item(16, 30).Value = "D"
$ws2.Cells.Item(16, 31).Value = "D"
$ws2.Cells.Item(16, 32).Value = "D"
$ws2.Cells.Item(16, 33).Value = "X"
$ws2.Cells.Item(16, 34).Value = "M"
$ws2.Cells.Item(16, 35).Value = "W"
$ws2.Cells.Item(17, 28).Value = "F"
$ws2.Cells.Item(17, 29).Value = "F"
$ws2.Cells.Item(17, 30).Value = "F"
$ws2.Cells.Item(17, 31).Value = "F"
$ws2.Cells.Item(17, 32).Value = "F"
$ws2.Cells.Item(17, 33).Value = "D"
$ws2.Cells.Item(17, 34).Value = "X"
$ws2.Cells.Item(17, 35).Value = "M"
$ws2.Cells.Item(17, 36).Value = "W"
$ws2.Cells.Item(18, 1).Value = "lw r4, 4(r4)"
$ws2.Cells.Item(18, 33).Value = "F"
$ws2.Cells.Item(18, 34).Value = "D"
$ws2.Cells.Item(18, 35).Value = "X"
$ws2.Cells.Item(18, 36).Value = "M"
$ws2.Cells.Item(18, 37).Value = "W"
$ws2.Cells.Item(19, 1).Value = "lw r5, 8(r5)"
$ws2.Cells.Item(19, 34).Value = "F"
$ws2.Cells.Item(19, 35).Value = "D"
$ws2.Cells.Item(19, 36).Value = "X"
$ws2.Cells.Item(19, 37).Value = "M"
$ws2.Cells.Item(19, 38).Value = "W"
$ws2.Cells.Item(20, 35).Value = "F"
$ws2.Cells.Item(20, 36).Value = "D"
$ws2.Cells.Item(20, 37).Value = "X"
$ws2.Cells.Item(20, 38).Value = "M"
$ws2.Cells.Item(20, 39).Value = "W"
$ws2.Cells.Item(21, 36).Value = "F"
$ws2.Cells.Item(21, 37).Value = "D"
$ws2.Cells.Item(21, 38).Value = "-"
$ws2.Cells.Item(21, 39).Value = "-"
$ws2.Cells.Item(21, 40).Value = "-"
$ws2.Cells.Item(22, 37).Value = "F"
$ws2.Cells.Item(22, 38).Value = "-"
$ws2.Cells.Item(22, 39).Value = "-"
$ws2.Cells.Item(22, 40).Value = "-"
$ws2.Cells.Item(22, 41).Value = "-"
$ws2.Cells.Item(23, 38).Value = "F"
$ws2.Cells.Item(23, 39).Value = "D"
$ws2.Cells.Item(23, 40).Value = "X"
$ws2.Cells.Item(23, 41).Value = "M"
$ws2.Cells.Item(23, 42).Value = "W"

# --- clear the cells that parta had but partb does not ---
$ws2.Cells.Item(17, 24).ClearContents()
$ws2.Cells.Item(17, 25).ClearContents()
$ws2.Cells.Item(17, 26).ClearContents()
$ws2.Cells.Item(17, 27).ClearContents()
$ws2.Cells.Item(18, 25).ClearContents()
$ws2.Cells.Item(18, 26).ClearContents()
$ws2.Cells.Item(18, 27).ClearContents()
$ws2.Cells.Item(18, 28).ClearContents()
$ws2.Cells.Item(18, 29).ClearContents()
$ws2.Cells.Item(19, 26).ClearContents()
$ws2.Cells.Item(19, 27).ClearContents()
$ws2.Cells.Item(19, 28).ClearContents()
$ws2.Cells.Item(19, 29).ClearContents()
$ws2.Cells.Item(19, 30).ClearContents()
$ws2.Cells.Item(20, 27).ClearContents()
$ws2.Cells.Item(20, 28).ClearContents()
$ws2.Cells.Item(20, 29).ClearContents()
$ws2.Cells.Item(20, 30).ClearContents()
$ws2.Cells.Item(20, 31).ClearContents()
$ws2.Cells.Item(21, 28).ClearContents()
$ws2.Cells.Item(21, 29).ClearContents()
$ws2.Cells.Item(21, 30).ClearContents()
$ws2.Cells.Item(21, 31).ClearContents()
$ws2.Cells.Item(21, 32).ClearContents()
$ws2.Cells.Item(22, 29).ClearContents()
$ws2.Cells.Item(22, 30).ClearContents()
$ws2.Cells.Item(22, 31).ClearContents()
$ws2.Cells.Item(22, 32).ClearContents()
$ws2.Cells.Item(22, 33).ClearContents()
$ws2.Cells.Item(23, 30).ClearContents()
$ws2.Cells.Item(23, 31).ClearContents()
$ws2.Cells.Item(23, 32).ClearContents()
$ws2.Cells.Item(23, 33).ClearContents()
$ws2.Cells.Item(23, 34).ClearContents()

# --- selections: parta no longer the active tab; partb is ---
$ws1.Range("X1").Select()
$ws2.Select()
$ws2.Range("AQ1").Select()
